$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '72.235.76'
$ws.Range('E2').Value = '  +3.80%  '
$ws.Range('D3').Value = '4.054.43'
$ws.Range('E3').Value = '  +3.35%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '521.79'
$ws.Range('E5').Value = '  -1.77%  '
$ws.Range('D6').Value = '148.56'
$ws.Range('E6').Value = '  +2.77%  '
$ws.Range('E7').Value = '  +1.51%  '
$ws.Range('E8').Value = '  +0.20%  '
$ws.Range('D9').Value = '0.739'
$ws.Range('E9').Value = '  +1.58%  '
$ws.Range('D10').Value = '0.175'
$ws.Range('E10').Value = '  +1.82%  '
$ws.Range('D11').Value = '0.0000338'
$ws.Range('E11').Value = '  +1.37%  '
$ws.Range('D12').Value = '47.06'
$ws.Range('E12').Value = '  +10.91%  '
$ws.Range('D13').Value = '10.81'
$ws.Range('E13').Value = '  +5.10%  '
$ws.Range('D14').Value = '4.705.75'
$ws.Range('E14').Value = '  +3.40%  '
$ws.Range('D15').Value = '4.090.57'
$ws.Range('E15').Value = '  +4.13%  '
$ws.Range('D16').Value = '21.42'
$ws.Range('E16').Value = '  +7.76%  '
$ws.Range('D17').Value = '14.29'
$ws.Range('E17').Value = '  +2.90%  '
$ws.Range('E18').Value = '  -1.14%  '
$ws.Range('D20').Value = '72.184.61'
$ws.Range('E20').Value = '  +3.85%  '
$ws.Range('D21').Value = '444.36'
$ws.Range('E21').Value = '  +1.77%  '
$ws.Range('B22').Value = 'Litecoin'
$ws.Range('C22').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D22').Value = '96.08'
$ws.Range('E22').Value = '  +9.15%  '
$ws.Range('B23').Value = 'ImmutableX'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D23').Value = '3.53'
$ws.Range('E23').Value = '  +5.16%  '
$ws.Range('D24').Value = '14.47'
$ws.Range('E24').Value = '  +0.41%  '
$ws.Range('D25').Value = '12.10'
$ws.Range('E25').Value = '  +4.69%  '
$ws.Range('D26').Value = '4.06'
$ws.Range('E26').Value = '  -1.61%  '
$ws.Range('D27').Value = '11.27'
$ws.Range('E27').Value = '  +5.13%  '
$ws.Range('D28').Value = '37.21'
$ws.Range('E28').Value = '  +2.23%  '
$ws.Range('D29').Value = '5.77'
$ws.Range('E29').Value = '  +2.10%  '
$ws.Range('D30').Value = '711.93'
$ws.Range('E30').Value = '  +2.10%  '
$ws.Range('E31').Value = '  +8.71%  '
$ws.Range('E32').Value = '  +2.24%  '
$ws.Range('E33').Value = '  +2.89%  '
$ws.Range('D34').Value = '6.94'
$ws.Range('E34').Value = '  +15.68%  '
$ws.Range('D35').Value = '67.57'
$ws.Range('E35').Value = '  -1.52%  '
$ws.Range('D36').Value = '0.0₃0911'
$ws.Range('E36').Value = '  +8.59%  '
$ws.Range('D37').Value = '0.443'
$ws.Range('E37').Value = '  -1.60%  '
$ws.Range('B38').Value = 'InjectiveProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D38').Value = '40.81'
$ws.Range('E38').Value = '  +1.13%  '
$ws.Range('B39').Value = 'ThetaToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D39').Value = '3.65'
$ws.Range('E39').Value = '  +23.49%  '
$ws.Range('E40').Value = '  +3.30%  '
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  -0.22%  '
$ws.Range('D43').Value = '0.0488'
$ws.Range('E43').Value = '  +0.41%  '
$ws.Range('E44').Value = '  +0.79%  '
$ws.Range('D45').Value = '2.81'
$ws.Range('E45').Value = '  +1.50%  '
$ws.Range('E46').Value = '  +5.03%  '
$ws.Range('E47').Value = '  +2.79%  '
$ws.Range('E48').Value = '  +2.26%  '
$ws.Range('D49').Value = '9.20'
$ws.Range('E49').Value = '  +7.82%  '
$ws.Range('E50').Value = '  +22.58%  '
$ws.Range('D51').Value = '3.34'
$ws.Range('E51').Value = '  +0.85%  '
